$wb = $excel.ActiveWorkbook

# --- Survey sheet: fix translation string column headers ---
# Column G header was "display.text" -> becomes the more specific "display.prompt.text"
# Column I header was "display.hint" -> becomes the more specific "display.hint.text"
$survey = $wb.Worksheets.Item("survey")
$survey.Range("G1").Value = "display.prompt.text"
$survey.Range("I1").Value = "display.hint.text"

# --- Settings sheet: bump form_version ---
$settings = $wb.Worksheets.Item("settings")
$settings.Range("B8").Value = 20170714

# --- View/selection bookkeeping to match the authored commit ---
# Survey sheet: scroll/selection moved to I2
$survey.Activate()
$excelWin = $excel.ActiveWindow
$excelWin.ScrollColumn = 5
$excelWin.ScrollRow = 1
$survey.Range("I2").Select()

# Properties sheet: keep its own selection (unchanged) but it is no longer the active tab
$properties = $wb.Worksheets.Item("properties")
$properties.Range("E6").Select()

# Settings sheet becomes the active tab, with its selection moved to C9
$settings.Activate()
$settings.Range("C9").Select()
